# Daily attendance processing - 2025-11-05 21:42:30
#
# The "Recorded By" column (G) lists the user(s)/process(es) that recorded
# each attendance session as a comma-separated string, e.g.
#   "dnasr281@gmail.com, System"
# For rows where "System" is the trailing entry in that list, move it to the
# front, e.g.
#   "System, dnasr281@gmail.com"
# Rows whose list does not end in "System" (e.g. a bare "dnasr281@gmail.com",
# or "admin@admin.com, dnasr281@gmail.com") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$col = 7  # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if ($value -eq "") { continue }

    $parts = $value -split ", "
    if ($parts.Count -lt 2) { continue }

    if ($parts[$parts.Count - 1] -eq "System") {
        $reordered = @("System") + $parts[0..($parts.Count - 2)]
        $cell.Value2 = $reordered -join ", "
    }
}
